$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for D2:E51 so numeric-looking strings (e.g. "1.00", "486.29")
# are not auto-converted to numbers by Excel, matching the inline-string cells in the diff.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "55.409.42"
$ws.Range("E2").Value = "  -4.04%  "
$ws.Range("D3").Value = "2.945.69"
$ws.Range("E3").Value = "  -6.71%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "486.29"
$ws.Range("E5").Value = "  -7.36%  "
$ws.Range("D6").Value = "129.88"
$ws.Range("E6").Value = "  -2.15%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "2.943.44"
$ws.Range("E8").Value = "  -6.77%  "
$ws.Range("E9").Value = "  -7.91%  "
$ws.Range("D10").Value = "7.01"
$ws.Range("E10").Value = "  -3.63%  "
$ws.Range("D11").Value = "0.0997"
$ws.Range("E11").Value = "  -10.07%  "
$ws.Range("E12").Value = "  -10.64%  "
$ws.Range("E13").Value = "  -0.65%  "
$ws.Range("D14").Value = "3.451.68"
$ws.Range("E14").Value = "  -6.70%  "
$ws.Range("D15").Value = "24.30"
$ws.Range("E15").Value = "  -6.08%  "
$ws.Range("D16").Value = "55.345.69"
$ws.Range("E16").Value = "  -4.19%  "
$ws.Range("D17").Value = "2.945.87"
$ws.Range("E17").Value = "  -6.74%  "
$ws.Range("E18").Value = "  -9.36%  "
$ws.Range("D19").Value = "5.61"
$ws.Range("E19").Value = "  -3.47%  "
$ws.Range("D20").Value = "11.95"
$ws.Range("E20").Value = "  -8.40%  "
$ws.Range("E21").Value = "  -8.17%  "
$ws.Range("D22").Value = "311.45"
$ws.Range("E22").Value = "  -9.94%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("E24").Value = "  -10.35%  "
$ws.Range("D25").Value = "59.76"
$ws.Range("E25").Value = "  -14.33%  "
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("D27").Value = "0.157"
$ws.Range("E27").Value = "  -5.70%  "
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("D29").Value = "0.0₃0842"
$ws.Range("E29").Value = "  -12.16%  "
$ws.Range("D30").Value = "6.52"
$ws.Range("E30").Value = "  -4.83%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "6.45"
$ws.Range("E31").Value = "  -6.81%  "
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "1.15"
$ws.Range("E32").Value = "  -4.55%  "
$ws.Range("E33").Value = "  -10.87%  "
$ws.Range("D34").Value = "19.20"
$ws.Range("E34").Value = "  -11.56%  "
$ws.Range("D35").Value = "147.31"
$ws.Range("E35").Value = "  -7.63%  "
$ws.Range("D36").Value = "4.36"
$ws.Range("E36").Value = "  -11.16%  "
$ws.Range("D37").Value = "5.62"
$ws.Range("E37").Value = "  -10.19%  "
$ws.Range("E38").Value = "  -9.50%  "
$ws.Range("D39").Value = "23.55"
$ws.Range("E39").Value = "  -9.53%  "
$ws.Range("E40").Value = "  -7.24%  "
$ws.Range("D41").Value = "2.975.57"
$ws.Range("E41").Value = "  -6.63%  "
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("D43").Value = "36.12"
$ws.Range("E43").Value = "  -11.21%  "
$ws.Range("D44").Value = "0.997"
$ws.Range("E44").Value = "  -7.60%  "
$ws.Range("D45").Value = "0.630"
$ws.Range("E45").Value = "  -9.79%  "
$ws.Range("E46").Value = "  -6.12%  "
$ws.Range("E47").Value = "  -10.99%  "
$ws.Range("D48").Value = "2.108.27"
$ws.Range("E48").Value = "  -7.07%  "
$ws.Range("E49").Value = "  -3.19%  "
$ws.Range("D50").Value = "18.86"
$ws.Range("E50").Value = "  -7.97%  "
$ws.Range("D51").Value = "5.53"
$ws.Range("E51").Value = "  -10.65%  "

# Restore default (unstyled) formatting on D2:E51 so the cell style index matches the original
# (these columns had no explicit style in the source workbook).
$ws.Range("D2:E51").ClearFormats()
